# Auto-generated Excel COM-interop script to apply diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # 展览
$ws.Range("F5").Value = 1318
$ws.Range("G5").Value = 78
$ws.Range("G6").Value = 70
$ws.Range("F8").Value = 928
$ws.Range("G8").Value = 19.9
$ws.Range("F9").Value = 738
$ws.Range("F10").Value = 210
$ws.Range("F15").Value = 3050
$ws.Range("F16").Value = 2676
$ws.Range("F18").Value = 32
$ws.Range("F23").Value = 5462
$ws.Range("F24").Value = 599
$ws.Range("F25").Value = 1001
$ws.Range("F27").Value = 68
$ws.Range("F28").Value = 381
$ws.Range("F29").Value = 1153
$ws.Range("F32").Value = 306
$ws = $wb.Worksheets.Item(2)  # 演出
$ws.Range("F3").Value = 1164
$ws.Range("F25").Value = 4006
$ws = $wb.Worksheets.Item(3)  # 本地生活
$ws.Range("F2").Value = 1790
$ws.Range("F5").Value = 2525
$ws.Range("F6").Value = 1084
$ws.Range("F9").Value = 1386
$ws.Range("F10").Value = 387
$ws = $wb.Worksheets.Item(4)  # 全部类型
$ws.Range("F2").Value = 1790
$ws.Range("F5").Value = 2525
$ws.Range("F9").Value = 1084
$ws.Range("F10").Value = 1386
$ws.Range("F11").Value = 387
$ws.Range("F14").Value = 1318
$ws.Range("G14").Value = 78
$ws.Range("G15").Value = 70
$ws.Range("F16").Value = 928
$ws.Range("G16").Value = 19.9
$ws.Range("F17").Value = 738
$ws.Range("F18").Value = 1164
$ws.Range("F19").Value = 210
$ws.Range("F23").Value = 3050
$ws.Range("F24").Value = 2676
$ws.Range("F25").Value = 32
$ws.Range("F31").Value = 5462
$ws.Range("F32").Value = 599
$ws.Range("F33").Value = 1001
$ws.Range("F36").Value = 68
$ws.Range("F37").Value = 381
$ws.Range("F50").Value = 306
